$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1) Remove the "Responsables" sheet entirely
# ------------------------------------------------------------------
$wsResp = $wb.Worksheets.Item("Responsables")
$wsResp.Delete()

# ------------------------------------------------------------------
# 2) OrdenEstandar: rename "Enchapado" -> "Barniz" and insert the new
#    processes "Stamping", "Cuño" and "Encapado" ahead of "Troquelado",
#    renumbering the Secuencia column below.
# ------------------------------------------------------------------
$wsOrden = $wb.Worksheets.Item("OrdenEstandar")

# Row 4 currently holds "Enchapado" -> becomes "Barniz"
$wsOrden.Range("B4").Value = "Barniz"

# Insert three new rows before the current row 6 ("Troquelado")
$wsOrden.Rows.Item(6).EntireRow.Insert()
$wsOrden.Rows.Item(6).EntireRow.Insert()
$wsOrden.Rows.Item(6).EntireRow.Insert()

$wsOrden.Range("A6").Value = 5
$wsOrden.Range("B6").Value = "Stamping"
$wsOrden.Range("A7").Value = 6
$wsOrden.Range("B7").Value = "Cuño"
$wsOrden.Range("A8").Value = 7
$wsOrden.Range("B8").Value = "Encapado"

# Renumber the remaining (shifted) rows: Troquelado, Descartonado, Ventana, Pegado
$wsOrden.Range("A9").Value = 8
$wsOrden.Range("A10").Value = 9
$wsOrden.Range("A11").Value = 10
$wsOrden.Range("A12").Value = 11

# ------------------------------------------------------------------
# 3) Maquinas: add "Troquelado / Manual 2" row, drop "Pegado / Manual"
#    row, and fix the "Enchapado" -> "Encapado" typo.
# ------------------------------------------------------------------
$wsMaq = $wb.Worksheets.Item("Maquinas")

# Insert a new row before current row 4 ("Pegado" / "Automática")
$wsMaq.Rows.Item(4).EntireRow.Insert()
$wsMaq.Range("A4").Value = "Troquelado"
$wsMaq.Range("B4").Value = "Manual 2"
$wsMaq.Range("C4").Value = 1000
$wsMaq.Range("D4").Value = 25
$wsMaq.Range("E4").Value = 10

# The old "Pegado / Manual / 1200 / 20 / 8" row is now at row 6; remove it
$wsMaq.Rows.Item(6).EntireRow.Delete()

# Fix "Enchapado" -> "Encapado"
$wsMaq.Range("A12").Value = "Encapado"

# ------------------------------------------------------------------
# 4) Jornada: give the (still empty) B4 cell a style, matching the
#    lightly-reformatted layout from the authoring session.
# ------------------------------------------------------------------
$wsJornada = $wb.Worksheets.Item("Jornada")
$wsJornada.Range("B4").Font.Bold = $false
